$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('New Account Creation', 'Passed'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Fail'),
    @('Bill Pay', 'Fail'),
    @('Login with username and password', 'Fail'),
    @('Login with username and password', 'Fail'),
    @('Loan Application', 'Fail'),
    @('Register with username and password', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Passed'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Update Contact Info', 'Passed'),
    @('Inter-Account Funds', 'Passed'),
    @('Loan Application', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Fail'),
    @('Inter-Account Funds', 'Fail'),
    @('Loan Application', 'Fail'),
    @('Login with username and password', 'Fail'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Fail'),
    @('Login with username and password', 'Fail'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Fail'),
    @('Login with username and password', 'Fail'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Fail'),
    @('New Account Creation', 'Passed'),
    @('Bill Pay', 'Passed'),
    @('Inter-Account Funds', 'Fail'),
    @('Loan Application', 'Passed'),
    @('Login with username and password', 'Passed'),
    @('Register with username and password', 'Fail')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A18:A19").Select()
